$wb = $excel.ActiveWorkbook

# --- Private sheet updates ---
$priv = $wb.Worksheets.Item("Private")
$priv.Range("B11").Value = "Color"
$priv.Range("E2").Value = "Raise"
$priv.Range("E11").Value = "Series E"
$priv.Range("F2").Value = "Amount"
$priv.Range("F11").Value = 482

# --- Main sheet updates ---
$main = $wb.Worksheets.Item("Main")
$main.Range("A3").Value = "x"
$main.Range("B3").Value = "Oracle"
$main.Range("B4").Value = "Epic"

$priv.Select()
$priv.Range("F11").Select()

$main.Select()
$main.Range("B5").Select()
